$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / "last updated" text (A1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Julio de 2020 a las 00:17"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 3824120
$ws.Range("C4").Value = 54108
$ws.Range("D4").Value = 1773406
$ws.Range("E4").Value = 1907946
$ws.Range("G4").Value = 704
$ws.Range("H4").Value = 142768

# --- Row 5: Brasil ---
$ws.Range("B5").Value = 2075124
$ws.Range("C5").Value = 26427
$ws.Range("E5").Value = 629614
$ws.Range("G5").Value = 803
$ws.Range("H5").Value = 78735

# --- Row 6: India ---
$ws.Range("B6").Value = 1077864
$ws.Range("C6").Value = 37407
$ws.Range("D6").Value = 677630
$ws.Range("E6").Value = 373406

# --- Row 9: Peru ---
$ws.Range("B9").Value = 349500
$ws.Range("C9").Value = 3963
$ws.Range("D9").Value = 238086
$ws.Range("E9").Value = 98416
$ws.Range("G9").Value = 199
$ws.Range("H9").Value = 12998

# --- Row 19: Alemania ---
$ws.Range("B19").Value = 202572
$ws.Range("C19").Value = 227
$ws.Range("E19").Value = 5910

# --- Row 24: Canada ---
$ws.Range("B24").Value = 109999
$ws.Range("C24").Value = 330
$ws.Range("D24").Value = 96912
$ws.Range("E24").Value = 4239

# --- Row 27: Egipto ---
$ws.Range("B27").Value = 87172
$ws.Range("C27").Value = 698
$ws.Range("D27").Value = 27868
$ws.Range("E27").Value = 55053
$ws.Range("G27").Value = 63
$ws.Range("H27").Value = 4251

# --- Row 44: Israel ---
$ws.Range("B44").Value = 49365
$ws.Range("C44").Value = 1906
$ws.Range("D44").Value = 21348
$ws.Range("E44").Value = 27616
$ws.Range("G44").Value = 9
$ws.Range("H44").Value = 401

# --- Row 49: Barein ---
$ws.Range("B49").Value = 36004
$ws.Range("C49").Value = 531
$ws.Range("D49").Value = 31765
$ws.Range("E49").Value = 4115

# --- Rows 66-70: Uzbekistan moves above Camerun, Costa de Marfil moves above Chequia ---
# Row 66 -> Uzbekistan (new data)
$ws.Range("A66").Value = "Uzbekistan"
$ws.Range("B66").Value = 16186
$ws.Range("C66").Value = 579
$ws.Range("D66").Value = 9127
$ws.Range("E66").Value = 6976
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 4
$ws.Range("H66").Value = 83

# Row 67 -> Camerun (previous row-66 data, unchanged)
$ws.Range("A67").Value = "Camerun"
$ws.Range("B67").Value = 16157
$ws.Range("C67").Value = 0
$ws.Range("D67").Value = 13728
$ws.Range("E67").Value = 2056
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 373

# Row 68 -> Costa de Marfil (new data)
$ws.Range("A68").Value = "Costa de Marfil"
$ws.Range("B68").Value = 13912
$ws.Range("C68").Value = 216
$ws.Range("D68").Value = 8000
$ws.Range("E68").Value = 5821
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 91

# Row 69 -> Chequia (previous row-68 data, unchanged)
$ws.Range("A69").Value = "Chequia"
$ws.Range("B69").Value = 13795
$ws.Range("C69").Value = 53
$ws.Range("D69").Value = 8725
$ws.Range("E69").Value = 4712
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 358

# Row 70 -> Corea del Sur (previous row-69 data, unchanged)
$ws.Range("A70").Value = "Corea del Sur"
$ws.Range("B70").Value = 13711
$ws.Range("C70").Value = 39
$ws.Range("D70").Value = 12519
$ws.Range("E70").Value = 898
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 294

# Row 71 (Costa de Marfil's old slot, now Dinamarca) is unchanged.

# --- Row 84: Bulgaria ---
$ws.Range("B84").Value = 8638
$ws.Range("C84").Value = 196
$ws.Range("D84").Value = 4081
$ws.Range("E84").Value = 4258
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 299

# --- Row 95: Mauritania ---
$ws.Range("B95").Value = 5813
$ws.Range("C95").Value = 103
$ws.Range("D95").Value = 3279
$ws.Range("E95").Value = 2381
$ws.Range("G95").Value = 2
$ws.Range("H95").Value = 153

# --- Row 128: Sierra Leona ---
$ws.Range("B128").Value = 1701
$ws.Range("C128").Value = 13
$ws.Range("D128").Value = 1237
$ws.Range("E128").Value = 399

# --- Row 133: Ruanda ---
$ws.Range("B133").Value = 1539
$ws.Range("C133").Value = 54
$ws.Range("D133").Value = 819
$ws.Range("E133").Value = 715
$ws.Range("G133").Value = 1
$ws.Range("H133").Value = 5

# --- Rows 150-151: Togo moves above Jamaica ---
# Row 150 -> Togo (new data)
$ws.Range("A150").Value = "Togo"
$ws.Range("B150").Value = 774
$ws.Range("C150").Value = 8
$ws.Range("D150").Value = 548
$ws.Range("E150").Value = 211
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 15

# Row 151 -> Jamaica (previous row-150 data, unchanged)
$ws.Range("A151").Value = "Jamaica"
$ws.Range("B151").Value = 768
$ws.Range("C151").Value = 3
$ws.Range("D151").Value = 678
$ws.Range("E151").Value = 80
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 10
